$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data block (rows 2-13) reflecting the updated NATMI LR-pair computation
# following Dr Hou advice: every sending cluster is now paired against both target
# clusters (FAPs and sCs) instead of only FAPs, and the underlying stats are recomputed.
$data = @(
    ,("ECs", "Gnai2", "Agtr2", "FAPs", 3, 1, 106.8785896666667, 320.635769, 0.1508748302900445, 0.1508748302900445, 3, 1, 0.9721403333333333, 2.916421, 0.4203533802573702, 0.4203533802573701, 103.9009877847499, 935.1088900627489, 0.06342074490817728, 0.06342074490817727)
    ,("ECs", "Gnai2", "Agtr2", "sCs", 3, 1, 106.8785896666667, 320.635769, 0.1508748302900445, 0.1508748302900445, 2, 0.6666666666666666, 1.340533666666667, 4.021601, 0.5796466197426299, 0.5796466197426299, 143.2743476940188, 1289.469129246169, 0.08745408538186726, 0.08745408538186726)
    ,("FAPs", "Gnai2", "Agtr2", "FAPs", 3, 1, 77.232923, 231.698769, 0.1090256166999485, 0.1090256166999485, 3, 1, 0.9721403333333333, 2.916421, 0.4203533802573702, 0.4203533802573701, 75.08123950952766, 675.731155585749, 0.04582928651446774, 0.04582928651446774)
    ,("FAPs", "Gnai2", "Agtr2", "sCs", 3, 1, 77.232923, 231.698769, 0.1090256166999485, 0.1090256166999485, 2, 0.6666666666666666, 1.340533666666667, 4.021601, 0.5796466197426299, 0.5796466197426299, 103.5333334565743, 931.8000011091691, 0.06319633018548076, 0.06319633018548078)
    ,("M1", "Gnai2", "Agtr2", "FAPs", 3, 1, 136.676337, 410.029011, 0.1929387280825172, 0.1929387280825172, 3, 1, 0.9721403333333333, 2.916421, 0.4203533802573702, 0.4203533802573701, 132.868579809959, 1195.817218289631, 0.08110244653204372, 0.08110244653204371)
    ,("M1", "Gnai2", "Agtr2", "sCs", 3, 1, 136.676337, 410.029011, 0.1929387280825172, 0.1929387280825172, 2, 0.6666666666666666, 1.340533666666667, 4.021601, 0.5796466197426299, 0.5796466197426299, 183.219231185179, 1648.973080666611, 0.1118362815504735, 0.1118362815504735)
    ,("M2", "Gnai2", "Agtr2", "FAPs", 3, 1, 143.4557186666667, 430.367156, 0.2025088212285795, 0.2025088212285795, 3, 1, 0.9721403333333333, 2.916421, 0.4203533802573702, 0.4203533802573701, 139.4590901631862, 1255.131811468676, 0.08512526753536889, 0.08512526753536888)
    ,("M2", "Gnai2", "Agtr2", "sCs", 3, 1, 143.4557186666667, 430.367156, 0.2025088212285795, 0.2025088212285795, 2, 0.6666666666666666, 1.340533666666667, 4.021601, 0.5796466197426299, 0.5796466197426299, 192.3072205485284, 1730.764984936756, 0.1173835536932106, 0.1173835536932106)
    ,("Neutro", "Gnai2", "Agtr2", "FAPs", 3, 1, 198.5982106666667, 595.794632, 0.2803505493821544, 0.2803505493821544, 3, 1, 0.9721403333333333, 2.916421, 0.4203533802573702, 0.4203533802573701, 193.0653307168969, 1737.587976452072, 0.1178463010897994, 0.1178463010897993)
    ,("Neutro", "Gnai2", "Agtr2", "sCs", 3, 1, 198.5982106666667, 595.794632, 0.2803505493821544, 0.2803505493821544, 2, 0.6666666666666666, 1.340533666666667, 4.021601, 0.5796466197426299, 0.5796466197426299, 266.2275875384258, 2396.048287845832, 0.162504248292355, 0.162504248292355)
    ,("sCs", "Gnai2", "Agtr2", "FAPs", 3, 1, 45.55066433333334, 136.651993, 0.06430145431675577, 0.06430145431675577, 3, 1, 0.9721403333333333, 2.916421, 0.4203533802573702, 0.4203533802573701, 44.28163800856144, 398.534742077053, 0.02702933367751316, 0.02702933367751315)
    ,("sCs", "Gnai2", "Agtr2", "sCs", 3, 1, 45.55066433333334, 136.651993, 0.06430145431675577, 0.06430145431675577, 2, 0.6666666666666666, 1.340533666666667, 4.021601, 0.5796466197426299, 0.5796466197426299, 61.0621990778659, 549.5597917007931, 0.03727212063924262, 0.03727212063924262)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
